$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 -- shifts the existing rows 17-24 down to 18-25
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with this week's data point (weekly Cilantro price update)
$ws.Cells.Item(17, 1).Value = 5
$ws.Cells.Item(17, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(17, 3).Value = 'Maule'
$ws.Cells.Item(17, 4).Value = 44736
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 6).Value = 100112040
$ws.Cells.Item(17, 7).Value = 'Cilantro'
$ws.Cells.Item(17, 8).Value = 'Sin especificar'
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 150
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 12).Value = 9000
$ws.Cells.Item(17, 13).Value = 9000
$ws.Cells.Item(17, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(17, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 16).Value = 250
$ws.Cells.Item(17, 17).Value = 36
$ws.Cells.Item(17, 18).Value = 'Hortaliza'
